# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price values that look like plain decimals are entered with a leading
# apostrophe so Excel keeps them as text (matching the sheet's original
# inline-string "Price" column) instead of auto-coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.754.18'
$ws.Range("E2").Value = '  -3.39%  '
$ws.Range("D3").Value = '2.554.58'
$ws.Range("E3").Value = '  -1.80%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'503.71"
$ws.Range("E5").Value = '  -3.75%  '
$ws.Range("D6").Value = "'141.52"
$ws.Range("E6").Value = '  -8.07%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = "'0.554"
$ws.Range("E8").Value = '  -5.53%  '
$ws.Range("D9").Value = '2.560.18'
$ws.Range("E9").Value = '  -1.75%  '
$ws.Range("D10").Value = "'6.20"
$ws.Range("E10").Value = '  -7.21%  '
$ws.Range("E11").Value = '  -4.32%  '
$ws.Range("E12").Value = '  -4.97%  '
$ws.Range("E13").Value = '  -1.02%  '
$ws.Range("D14").Value = '3.004.76'
$ws.Range("E14").Value = '  -1.67%  '
$ws.Range("D15").Value = '58.799.01'
$ws.Range("E15").Value = '  -3.33%  '
$ws.Range("D16").Value = "'20.45"
$ws.Range("E16").Value = '  -5.24%  '
$ws.Range("E17").Value = '  -4.91%  '
$ws.Range("D18").Value = '2.557.42'
$ws.Range("E18").Value = '  -1.73%  '
$ws.Range("E19").Value = '  -5.37%  '
$ws.Range("D20").Value = "'331.52"
$ws.Range("E20").Value = '  -6.53%  '
$ws.Range("E21").Value = '  -5.13%  '
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").Value = "'5.92"
$ws.Range("E23").Value = '  -4.31%  '
$ws.Range("D24").Value = "'59.64"
$ws.Range("E24").Value = '  -2.40%  '
$ws.Range("E25").Value = '  -5.11%  '
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("D27").Value = "'0.159"
$ws.Range("E27").Value = '  -4.18%  '
$ws.Range("E28").Value = '  -7.99%  '
$ws.Range("E29").Value = '  -7.12%  '
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").Value = "'148.79"
$ws.Range("E31").Value = '  -0.96%  '
$ws.Range("D32").Value = "'18.49"
$ws.Range("E32").Value = '  -4.55%  '
$ws.Range("E33").Value = '  -4.26%  '
$ws.Range("E34").Value = '  -8.06%  '
$ws.Range("D35").Value = "'3.86"
$ws.Range("E35").Value = '  -7.35%  '
$ws.Range("D36").Value = "'0.876"
$ws.Range("E36").Value = '  -4.97%  '
$ws.Range("D37").Value = "'1.09"
$ws.Range("E37").Value = '  -8.24%  '
$ws.Range("E38").Value = '  -1.60%  '
$ws.Range("E39").Value = '  -9.48%  '
$ws.Range("D40").Value = "'285.92"
$ws.Range("E40").Value = '  -4.07%  '
$ws.Range("E41").Value = '  -7.57%  '
$ws.Range("E42").Value = '  -7.56%  '
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("D44").Value = "'0.0980"
$ws.Range("E44").Value = '  -3.34%  '
$ws.Range("E45").Value = '  -2.86%  '
$ws.Range("E46").Value = '  -5.60%  '
$ws.Range("D47").Value = "'10.33"
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("E48").Value = '  -5.16%  '
$ws.Range("E49").Value = '  -5.36%  '
$ws.Range("D50").Value = "'4.49"
$ws.Range("E50").Value = '  -8.30%  '
$ws.Range("D51").Value = '1.887.81'
$ws.Range("E51").Value = '  -4.13%  '
